$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are stored as text in this sheet (inline strings).
# For values that look like plain numbers, force text formatting first
# so Excel does not auto-convert them to a numeric type, then restore
# the default "Normal" style so no stray formatting is introduced.

$ws.Range('D2').Value = '26.152.65'
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('D3').Value = '1.653.91'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('E6').Value = '  +1.33%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06323'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.42'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07809'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.516'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').Value = '1.668.20'
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('D14').Value = '1.883.25'
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').Value = '0.0₅8168'
$ws.Range('E16').Value = '  +1.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.33'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('D18').Value = '26.139.57'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.597'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.08'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.997'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('E24').Value = '  +0.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.45'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1224'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.201'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.68%  '
$ws.Range('E28').Value = '  -0.78%  '
$ws.Range('E29').Value = '  +3.12%  '
$ws.Range('E30').Value = '  -3.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.273'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.543'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.261'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.588'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.804'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.01%  '
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9478'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5717'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01606'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8513'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.59%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.797'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.005'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.31%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.038.84'
$ws.Range('E43').Value = '  +3.86%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '103.82'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.34%  '
$ws.Range('D45').Value = '1.795.82'
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.34%  '
$ws.Range('D47').Value = '0.0₈105'
$ws.Range('E47').Value = '  -1.06%  '
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4356'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05154'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.830'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.18%  '
